# Apply updated profit figures to the Sheets workbook (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6000
$ws.Range("J51").Value = 6000
$ws.Range("L51").Value = 6000
$ws.Range("N51").Value = -6968
$ws.Range("H70").Value = 7006.3125
$ws.Range("I70").Value = 1750
$ws.Range("J70").Value = 7757.2144
$ws.Range("K70").Value = 5250
$ws.Range("L70").Value = 23271.6432
$ws.Range("M70").Value = -4980
$ws.Range("N70").Value = -23811.6432
$ws.Range("H73").Value = 7006.3125
$ws.Range("I73").Value = 1750
$ws.Range("J73").Value = 7757.2144
$ws.Range("K73").Value = 5250
$ws.Range("L73").Value = 23271.6432
$ws.Range("M73").Value = -4314
$ws.Range("N73").Value = -25143.6432
$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 4500
$ws.Range("K74").Value = 4500
$ws.Range("M74").Value = -3564
$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 4500
$ws.Range("K77").Value = 22500
$ws.Range("M77").Value = -17820
$ws.Range("H98").Value = 587.6
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H107").Value = 1500
$ws.Range("I107").Value = 1500
$ws.Range("K107").Value = 1500
$ws.Range("M107").Value = 420
$ws.Range("H113").Value = 2665
$ws.Range("I113").Value = 2997.5
$ws.Range("K113").Value = 2997.5
$ws.Range("M113").Value = 256.5
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 587.6
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H135").Value = 2813.8333
$ws.Range("I135").Value = 2664.3333
$ws.Range("K135").Value = 23978.9997
$ws.Range("M135").Value = -21443.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 42499.5
$ws.Range("J125").Value = 42499.5
$ws.Range("L125").Value = 42499.5
$ws.Range("N125").Value = -52339.5
$ws.Range("H132").Value = 1269.7646
$ws.Range("I132").Value = 1346.6154
$ws.Range("K132").Value = 4039.8462
$ws.Range("M132").Value = -1509.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8842.571
$ws.Range("J20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("N20").Value = -2494
$ws.Range("H22").Value = 627.5333000000001
$ws.Range("I22").Value = 627.5333000000001
$ws.Range("K22").Value = 627.5333000000001
$ws.Range("M22").Value = -454.5333000000001
$ws.Range("H99").Value = 3117.6206
$ws.Range("I99").Value = 3102.0715
$ws.Range("J99").Value = 3132.1333
$ws.Range("K99").Value = 3102.0715
$ws.Range("L99").Value = 3132.1333
$ws.Range("M99").Value = -1604.0715
$ws.Range("N99").Value = -6128.1333
$ws.Range("H107").Value = 709.1177
$ws.Range("I107").Value = 709.1177
$ws.Range("K107").Value = 709.1177
$ws.Range("M107").Value = 1210.8823
$ws.Range("H134").Value = 2061.4814
$ws.Range("I134").Value = 1784.6364
$ws.Range("K134").Value = 5353.9092
$ws.Range("M134").Value = -2818.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15120.375
$ws.Range("I99").Value = 12348.9
$ws.Range("K99").Value = 12348.9
$ws.Range("M99").Value = -10850.9
$ws.Range("H126").Value = 15120.375
$ws.Range("I126").Value = 12348.9
$ws.Range("K126").Value = 37046.7
$ws.Range("M126").Value = -34576.7
$ws.Range("H134").Value = 3355
$ws.Range("I134").Value = 2492
$ws.Range("K134").Value = 7476
$ws.Range("M134").Value = -4941

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 188.6
$ws.Range("I8").Value = 188.6
$ws.Range("K8").Value = 565.8
$ws.Range("M8").Value = -426.8
$ws.Range("H122").Value = 633.55554
$ws.Range("I122").Value = 549.5
$ws.Range("J122").Value = 700.8
$ws.Range("K122").Value = 4945.5
$ws.Range("L122").Value = 6307.2
$ws.Range("M122").Value = -2495.5
$ws.Range("N122").Value = -11207.2
$ws.Range("H129").Value = 2107
$ws.Range("I129").Value = 690.6667
$ws.Range("K129").Value = 2072.0001
$ws.Range("M129").Value = 2927.9999
$ws.Range("H131").Value = 989.25
$ws.Range("I131").Value = 841.3333
$ws.Range("J131").Value = 1433
$ws.Range("K131").Value = 2523.9999
$ws.Range("L131").Value = 4299
$ws.Range("M131").Value = 2516.0001
$ws.Range("N131").Value = -14379
$ws.Range("H132").Value = 1500
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H141").Value = 11030
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 17333.334
$ws.Range("J95").Value = 17333.334
$ws.Range("L95").Value = 17333.334
$ws.Range("N95").Value = -22825.334
$ws.Range("H102").Value = 3558.7646
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3824.25
$ws.Range("I16").Value = 4081.182
$ws.Range("J16").Value = 998
$ws.Range("K16").Value = 4081.182
$ws.Range("L16").Value = 998
$ws.Range("M16").Value = -3911.182
$ws.Range("N16").Value = -1338
$ws.Range("H19").Value = 33000
$ws.Range("I19").Value = 33000
$ws.Range("K19").Value = 33000
$ws.Range("M19").Value = -32830
$ws.Range("H61").Value = 12996
$ws.Range("J61").Value = 9987
$ws.Range("L61").Value = 9987
$ws.Range("N61").Value = -10391
$ws.Range("H113").Value = 12996
$ws.Range("J113").Value = 9987
$ws.Range("L113").Value = 9987
$ws.Range("N113").Value = -14327
$ws.Range("H122").Value = 7105.2
$ws.Range("I122").Value = 7444.9287
$ws.Range("K122").Value = 22334.7861
$ws.Range("M122").Value = -19884.7861
$ws.Range("H136").Value = 3665.6667
$ws.Range("I136").Value = 3199
$ws.Range("J136").Value = 5999
$ws.Range("K136").Value = 9597
$ws.Range("L136").Value = 17997
$ws.Range("M136").Value = -7047
$ws.Range("N136").Value = -23097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5056.857
$ws.Range("I81").Value = 6079.6
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 12159.2
$ws.Range("L81").Value = 5000
$ws.Range("M81").Value = -11098.2
$ws.Range("N81").Value = -7122
$ws.Range("H84").Value = 5056.857
$ws.Range("I84").Value = 6079.6
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 60796
$ws.Range("L84").Value = 25000
$ws.Range("M84").Value = -55492
$ws.Range("N84").Value = -35608
$ws.Range("H122").Value = 2962.5
$ws.Range("I122").Value = 2800
$ws.Range("J122").Value = 3125
$ws.Range("K122").Value = 8400
$ws.Range("L122").Value = 9375
$ws.Range("M122").Value = -5950
$ws.Range("N122").Value = -14275
$ws.Range("H136").Value = 1503.6511
$ws.Range("J136").Value = 4057.4
$ws.Range("L136").Value = 12172.2
$ws.Range("N136").Value = -17272.2
$ws.Range("H141").Value = 87997
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 87997
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 87997
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -98357

